$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 4938.048
$ws.Range("I9").Value = 121.42857
$ws.Range("J9").Value = 14571.286
$ws.Range("K9").Value = 121.42857
$ws.Range("L9").Value = 14571.286
$ws.Range("M9").Value = 47.57143000000001
$ws.Range("N9").Value = -14909.286

$ws.Range("H18").Value = 1547.7142
$ws.Range("I18").Value = 1770.909
$ws.Range("K18").Value = 1770.909
$ws.Range("M18").Value = -1486.909

$ws.Range("H46").Value = 849999.5
$ws.Range("I46").Value = 849999.5
$ws.Range("K46").Value = 2549998.5
$ws.Range("M46").Value = -2549879.5

$ws.Range("H60").Value = 849999.5
$ws.Range("I60").Value = 849999.5
$ws.Range("K60").Value = 2549998.5
$ws.Range("M60").Value = -2549514.5

$ws.Range("H113").Value = 8937.5
$ws.Range("I113").Value = 7357.143
$ws.Range("K113").Value = 7357.143
$ws.Range("M113").Value = -4103.143

$ws.Range("H124").Value = 96519.664
$ws.Range("J124").Value = 96519.664
$ws.Range("L124").Value = 96519.664
$ws.Range("N124").Value = -106339.664

$ws.Range("H135").Value = 1035.1904
$ws.Range("I135").Value = 1092.3529
$ws.Range("K135").Value = 9831.176100000001
$ws.Range("M135").Value = -7296.176100000001

$ws.Range("H137").Value = 1056.8096
$ws.Range("I137").Value = 946
$ws.Range("K137").Value = 2838
$ws.Range("M137").Value = -288

$ws.Range("H138").Value = 2092.1904
$ws.Range("I138").Value = 1333.2727
$ws.Range("J138").Value = 2927
$ws.Range("K138").Value = 3999.8181
$ws.Range("L138").Value = 8781
$ws.Range("M138").Value = 1140.1819
$ws.Range("N138").Value = -19061

$ws.Range("H141").Value = 5311
$ws.Range("I141").Value = 5814.9414
$ws.Range("K141").Value = 17444.8242
$ws.Range("M141").Value = -12264.8242

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 33833.668
$ws.Range("J34").Value = 34501
$ws.Range("L34").Value = 34501
$ws.Range("N34").Value = -35043

$ws.Range("H45").Value = 35715504
$ws.Range("I45").Value = 45455612
$ws.Range("K45").Value = 45455612
$ws.Range("M45").Value = -45455235

$ws.Range("H61").Value = 15628421
$ws.Range("J61").Value = 4688.84
$ws.Range("L61").Value = 4688.84
$ws.Range("N61").Value = -5112.84

$ws.Range("H97").Value = 2767.2632
$ws.Range("I97").Value = 775.6667
$ws.Range("J97").Value = 4559.7
$ws.Range("K97").Value = 775.6667
$ws.Range("L97").Value = 4559.7
$ws.Range("M97").Value = -279.6667
$ws.Range("N97").Value = -5551.7

$ws.Range("H110").Value = 1690.6
$ws.Range("I110").Value = 1688.4
$ws.Range("K110").Value = 1688.4
$ws.Range("M110").Value = 356.5999999999999

$ws.Range("H136").Value = 15628421
$ws.Range("J136").Value = 4688.84
$ws.Range("L136").Value = 14066.52
$ws.Range("N136").Value = -19166.52

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H99").Value = 2330.5557
$ws.Range("I99").Value = 2276.1538
$ws.Range("J99").Value = 2472
$ws.Range("K99").Value = 2276.1538
$ws.Range("L99").Value = 2472
$ws.Range("M99").Value = -778.1538
$ws.Range("N99").Value = -5468

$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("M104").ClearContents()

$ws.Range("H107").Value = 11828.238
$ws.Range("I107").Value = 14786.5625
$ws.Range("J107").Value = 2361.6
$ws.Range("K107").Value = 14786.5625
$ws.Range("L107").Value = 2361.6
$ws.Range("M107").Value = -12866.5625
$ws.Range("N107").Value = -6201.6

$ws.Range("H132").Value = 99146.336
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws.Range("H138").Value = 78948.164
$ws.Range("J138").Value = 78948.164
$ws.Range("L138").Value = 78948.164
$ws.Range("N138").Value = -89228.164

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4011
$ws.Range("I16").Value = 4011
$ws.Range("K16").Value = 4011
$ws.Range("M16").Value = -3724

$ws.Range("H58").Value = 64104700
$ws.Range("I58").Value = 83335310
$ws.Range("K58").Value = 83335310
$ws.Range("M58").Value = -83335107

$ws.Range("H105").Value = 2398.889
$ws.Range("I105").Value = 2398.889
$ws.Range("K105").Value = 2398.889
$ws.Range("M105").Value = -651.8890000000001

$ws.Range("H113").Value = 4011
$ws.Range("I113").Value = 4011
$ws.Range("K113").Value = 4011
$ws.Range("M113").Value = -1841

$ws.Range("H134").Value = 3342.2
$ws.Range("I134").Value = 3247.0833
$ws.Range("K134").Value = 9741.249899999999
$ws.Range("M134").Value = -7206.249899999999

$ws.Range("H136").Value = 64104700
$ws.Range("I136").Value = 83335310
$ws.Range("K136").Value = 250005930
$ws.Range("M136").Value = -250003380

$ws.Range("H138").Value = 79270.36
$ws.Range("J138").Value = 78675.766
$ws.Range("L138").Value = 78675.766
$ws.Range("N138").Value = -88955.766

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1253.3334
$ws.Range("J12").Value = 1222.3125
$ws.Range("L12").Value = 3666.9375
$ws.Range("N12").Value = -4012.9375

$ws.Range("H51").Value = 10073.429
$ws.Range("I51").Value = 252
$ws.Range("J51").Value = 14002
$ws.Range("K51").Value = 756
$ws.Range("L51").Value = 42006
$ws.Range("M51").Value = -296
$ws.Range("N51").Value = -42926

$ws.Range("H64").Value = 7686
$ws.Range("I64").Value = 6672
$ws.Range("J64").Value = 8700
$ws.Range("K64").Value = 20016
$ws.Range("L64").Value = 26100
$ws.Range("M64").Value = -19746
$ws.Range("N64").Value = -26640

$ws.Range("H67").Value = 7686
$ws.Range("I67").Value = 6672
$ws.Range("J67").Value = 8700
$ws.Range("K67").Value = 20016
$ws.Range("L67").Value = 26100
$ws.Range("M67").Value = -19080
$ws.Range("N67").Value = -27972

$ws.Range("H113").Value = 142860830
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 142860830
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 428582490
$ws.Range("N113").Value = -428586830
$ws.Range("M113").ClearContents()

$ws.Range("H126").Value = 12736.385
$ws.Range("I126").Value = 3898
$ws.Range("J126").Value = 14343.363
$ws.Range("K126").Value = 11694
$ws.Range("L126").Value = 43030.089
$ws.Range("M126").Value = -6754
$ws.Range("N126").Value = -52910.089

$ws.Range("H140").Value = 2037.4546
$ws.Range("I140").Value = 1602
$ws.Range("K140").Value = 4806
$ws.Range("M140").Value = 374

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 17511.5
$ws.Range("J38").Value = 17511.5
$ws.Range("L38").Value = 17511.5
$ws.Range("N38").Value = -18437.5

$ws.Range("H135").Value = 93720.25999999999
$ws.Range("J135").Value = 93720.25999999999
$ws.Range("L135").Value = 93720.25999999999
$ws.Range("N135").Value = -103860.26

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 697.6129
$ws.Range("I16").Value = 745.5714
$ws.Range("K16").Value = 745.5714
$ws.Range("M16").Value = -575.5714

$ws.Range("H40").Value = 3542.4092
$ws.Range("I40").Value = 3915.611
$ws.Range("J40").Value = 1863
$ws.Range("K40").Value = 3915.611
$ws.Range("L40").Value = 1863
$ws.Range("M40").Value = -3779.611
$ws.Range("N40").Value = -2135

$ws.Range("H46").Value = 4357.5
$ws.Range("I46").Value = 3830.5
$ws.Range("J46").Value = 5148
$ws.Range("K46").Value = 3830.5
$ws.Range("L46").Value = 5148
$ws.Range("M46").Value = -3642.5
$ws.Range("N46").Value = -5524

$ws.Range("H132").Value = 3796
$ws.Range("I132").Value = 3170.0625
$ws.Range("K132").Value = 9510.1875
$ws.Range("M132").Value = -6980.1875

$ws.Range("H136").Value = 7814384
$ws.Range("I136").Value = 1836.25
$ws.Range("J136").Value = 62502216
$ws.Range("K136").Value = 5508.75
$ws.Range("L136").Value = 187506648
$ws.Range("M136").Value = -2958.75
$ws.Range("N136").Value = -187511748

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1158.2413
$ws.Range("I113").Value = 1107.7894
$ws.Range("J113").Value = 1254.1
$ws.Range("K113").Value = 3323.3682
$ws.Range("L113").Value = 3762.3
$ws.Range("M113").Value = -1153.3682
$ws.Range("N113").Value = -8102.299999999999

$ws.Range("H122").Value = 2425.6936
$ws.Range("I122").Value = 2097.2856
$ws.Range("K122").Value = 6291.8568
$ws.Range("M122").Value = -3841.8568

$ws.Range("H132").Value = 2572.0435
$ws.Range("I132").Value = 2550.7368
$ws.Range("K132").Value = 7652.2104
$ws.Range("M132").Value = -5122.2104

$ws.Range("H136").Value = 29916936
$ws.Range("I136").Value = 2483672.2
$ws.Range("J136").Value = 200003170
$ws.Range("K136").Value = 7451016.600000001
$ws.Range("L136").Value = 600009510
$ws.Range("M136").Value = -7448466.600000001
$ws.Range("N136").Value = -600014610
